$d = $word.ActiveDocument

# Locate the paragraph that currently ends the "Main.cpp" bullet-list section:
# "Run find queries on another set of predefined characters, including one that isn't in the list"
$count = $d.Paragraphs.Count
$targetIdx = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Run find queries on another set of predefined characters*") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -eq 0) {
    throw "Could not locate anchor paragraph"
}

# New LList.h bullet-list items to insert after the anchor paragraph, in order.
# Level is 0-based, matching w:ilvl in the OOXML. Items with two runs use the
# "Runs" key (array of run texts) instead of a single "Text".
$items = @(
    @{ Level = 0; Text = "LList.h" },
    @{ Level = 1; Text = "transpose()" },
    @{ Level = 2; Text = "Swaps the accessed item with the item immediately preceding it" },
    @{ Level = 1; Text = "reorderCount()" },
    @{ Level = 2; Text = "Handler for reordering the list by count" },
    @{ Level = 2; Text = "Increments the count variable" },
    @{ Level = 1; Text = "reorderCountHelp()" },
    @{ Level = 2; Runs = @("Used to ", "remove and insert the accessed item where its count variable is >= the element preceding it") },
    @{ Level = 1; Text = "moveToFront()" },
    @{ Level = 2; Text = "Removes the accessed element and inserts it at the front of the list" }
)

function Escape-Xml([string]$s) {
    $s = $s -replace "&", "&amp;"
    $s = $s -replace "<", "&lt;"
    $s = $s -replace ">", "&gt;"
    return $s
}

$curIdx = $targetIdx
foreach ($item in $items) {
    $curPara = $d.Paragraphs.Item($curIdx)
    $curRange = $curPara.Range
    $curRange.InsertParagraphAfter()

    $curIdx = $curIdx + 1
    $newPara = $d.Paragraphs.Item($curIdx)
    $newRange = $newPara.Range

    # Word's ListFormat.ListLevelNumber is 1-based; w:ilvl in OOXML is 0-based.
    $newRange.ListFormat.ListLevelNumber = $item.Level + 1

    if ($item.ContainsKey("Runs")) {
        $runsXml = ""
        foreach ($runText in $item.Runs) {
            $escaped = Escape-Xml $runText
            if ($runText -match "^\s" -or $runText -match "\s$") {
                $runsXml = $runsXml + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
            } else {
                $runsXml = $runsXml + '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/></w:rPr><w:t>' + $escaped + '</w:t></w:r>'
            }
        }

        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
          '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body><w:p>' +
          '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="' + $item.Level + '"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/></w:rPr></w:pPr>' +
          $runsXml +
          '</w:p></w:body></w:document>' +
          '</pkg:xmlData></pkg:part></pkg:package>'

        $null = $newRange.InsertXML($xml)
    } else {
        $newRange.Text = $item.Text
    }
}
